$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2..44: RowNum, open_price(D), close_price(E), high_price(F), low_price(G), shares_outstanding(H)
$data = @(
    @(2, 55.5, 57.38999938964844, 59.9900016784668, 48.5, 42689006),
    @(3, 73.33999633789062, 90.20999908447266, 91.83999633789062, 70.38999938964844, 42689006),
    @(4, 116, 95.62999725341795, 119.9899978637695, 71, 42689006),
    @(5, 49.86999893188477, 62.81000137329102, 64.26999664306641, 48.7599983215332, 42689006),
    @(6, 39, 46.40000152587891, 47.15999984741211, 33.38999938964844, 42689006),
    @(7, 41.02999877929688, 41.34000015258789, 41.5, 35.2599983215332, 42689006),
    @(8, 57.95000076293945, 72.04000091552734, 72.98999786376953, 55.75, 42689006),
    @(9, 61.22000122070312, 61.52000045776367, 65.77999877929688, 53.70000076293945, 42689006),
    @(10, 50, 58.95000076293945, 60.52000045776367, 49.18000030517578, 42689006),
    @(11, 56.20000076293945, 58.54000091552734, 65.38999938964844, 55.77999877929688, 42689006),
    @(12, 50.22000122070312, 54.40000152587891, 54.93000030517578, 49.22000122070312, 42689006),
    @(13, 56.59999847412109, 54.27000045776367, 58.54999923706055, 52.45000076293945, 42689006),
    @(14, 49.91999816894531, 48.08000183105469, 51.08000183105469, 46.09000015258789, 42689006),
    @(15, 46.79999923706055, 48.86000061035156, 55.5, 46, 42689006),
    @(16, 39.04000091552734, 38.31000137329102, 40.58000183105469, 34.33000183105469, 42689006),
    @(17, 34.95000076293945, 39.9900016784668, 40.65000152587891, 30, 42689006),
    @(18, 37.97000122070312, 40.36999893188477, 42.20000076293945, 37.16999816894531, 42689006),
    @(19, 50.45000076293945, 37.90999984741211, 51.36000061035156, 36.68000030517578, 42689006),
    @(20, 50.04000091552734, 55.86999893188477, 56.90000152587891, 43.0099983215332, 42689006),
    @(21, 52.81999969482422, 54.7400016784668, 61.75, 52.29999923706055, 42689006),
    @(22, 59.40000152587891, 59.45000076293945, 73.58999633789062, 55.22000122070312, 42689006),
    @(23, 50.77999877929688, 56.7400016784668, 58.9900016784668, 47.90999984741211, 42689006),
    @(24, 45.41999816894531, 52.70000076293945, 53.54999923706055, 45.09999847412109, 42689006),
    @(25, 55.45999908447266, 78.12999725341797, 79.83000183105469, 53.65999984741211, 42689006),
    @(26, 97.86000061035156, 112.4899978637695, 128.9299926757812, 97.73000335693359, 42689006),
    @(27, 98.72000122070312, 100.379997253418, 102.5800018310547, 82.58999633789062, 42689006),
    @(28, 99.36000061035156, 103.5699996948242, 106.4300003051758, 89.30000305175781, 42689006),
    @(29, 186.9400024414062, 179.5200042724609, 207, 170.7799987792969, 42689006),
    @(30, 142.9400024414062, 139.7100067138672, 156.1000061035156, 118.5800018310547, 42689006),
    @(31, 82.08000183105469, 85.22000122070312, 90.37999725341795, 65.59999847412109, 42689006),
    @(32, 85.93000030517578, 67.87999725341797, 93.62999725341795, 67.31999969482422, 42689006),
    @(33, 56.36999893188477, 74.19999694824219, 76.59999847412109, 51.29000091552734, 42689006),
    @(34, 90, 94.30999755859376, 99.86000061035156, 89.33000183105469, 42689006),
    @(35, 61.84999847412109, 72.31999969482422, 83.16000366210938, 61.2400016784668, 42689006),
    @(36, 82.41999816894531, 62.15000152587891, 84.72000122070312, 57.20000076293945, 42689006),
    @(37, 44.7400016784668, 58.70999908447266, 58.81000137329102, 43.59000015258789, 42689006),
    @(38, 53.13000106811523, 55.84999847412109, 65.38999938964844, 50.84000015258789, 42689006),
    @(39, 45.18000030517578, 58.2599983215332, 59.9900016784668, 44.52000045776367, 42689006),
    @(40, 51.59999847412109, 59.70000076293945, 63.97999954223633, 39.68999862670898, 42689006),
    @(41, 56.54000091552734, 71.55000305175781, 81.31999969482422, 56.2400016784668, 42689006),
    @(42, 73.77999877929688, 61.43000030517578, 84, 60.04999923706055, 42689006),
    @(43, 48.95000076293945, 52.63999938964844, 66.30999755859375, 48.09000015258789, 42689006),
    @(44, 64, 82.48000335693359, 93.44000244140624, 62.61999893188477, 42689006)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]   # D: open_price
    $ws.Cells.Item($r, 5).Value = $row[2]   # E: close_price
    $ws.Cells.Item($r, 6).Value = $row[3]   # F: high_price
    $ws.Cells.Item($r, 7).Value = $row[4]   # G: low_price
    $ws.Cells.Item($r, 8).Value = $row[5]   # H: shares_outstanding
    $ws.Cells.Item($r, 9).Value = "AMBA"    # I: fixed_ticker -> now always AMBA
}
